$wb = $excel.ActiveWorkbook

# --- Home sheet: selection moves from column C to column D, and column D
#     (which already holds the calories_goal figures later copied into the
#     Calories sheet) gets an explicit best-fit-like width.
$homeWs = $wb.Worksheets.Item("Home")
$homeWs.Columns.Item(4).ColumnWidth = 11.1
$homeWs.Range("D1:D1048576").Select()

# --- Calories sheet: drop the "required_intake" column and replace the
#     sample row with the real calories_goal data pulled from the Home
#     sheet (one row per user).
$calWs = $wb.Worksheets.Item("Calories")

# Clear out column C ("required_intake") entirely - header + data.
$calWs.Range("C1:C11").ClearContents()

$values = @(
    @(2000, 2000),
    @(1800, 2000),
    @(2200, 2000),
    @(2000, 2000),
    @(1900, 2000),
    @(2000, 3500),
    @(2500, 4000),
    @(2300, 3500),
    @(2200, 3000),
    @(3000, 2500)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $calWs.Cells.Item($row, 1).Value = $values[$i][0]
    $calWs.Cells.Item($row, 2).Value = $values[$i][1]
}

# Preserve the centered look on the two rows that originally carried it.
$calWs.Range("A2").HorizontalAlignment = -4108
$calWs.Range("A7").HorizontalAlignment = -4108

$calWs.Activate()
$calWs.Range("D7").Select()
